# Updated symbol list on Sat Feb 11 20:58:07 UTC 2023 with GitHub Actions
# Applies per-cell edits mirroring a day's worth of coinranking.com scrape refresh:
# price (D) / 1h volume change (E) updates for most rows, plus a swap of the
# BTSEToken / MXToken rows (7 <-> 8) including their coin name + link cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (matches the workbook's original
# inline-string cells) without leaving a stray NumberFormat on the cell -
# Excel would otherwise auto-coerce numeric-looking / percent-looking
# strings like "308.99" or "0.34%" into real numbers.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}


# Row 2 - BNB
Set-TextValue $ws.Range("D2") "308.99"
Set-TextValue $ws.Range("E2") "0.34%"

# Row 3 - OKB
Set-TextValue $ws.Range("D3") "40.94"
Set-TextValue $ws.Range("E3") "0.26%"

# Row 4 - HuobiToken
Set-TextValue $ws.Range("D4") "5.119"
Set-TextValue $ws.Range("E4") "1.55%"

# Row 5 - Cronos
Set-TextValue $ws.Range("D5") "0.07626"
Set-TextValue $ws.Range("E5") "-0.12%"

# Row 6 - FTXToken
Set-TextValue $ws.Range("D6") "1.605"
Set-TextValue $ws.Range("E6") "-0.79%"

# Row 7 - BTSEToken -> MXToken (row swap)
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D7") "0.9094"
Set-TextValue $ws.Range("E7") "-0.03%"

# Row 8 - MXToken -> BTSEToken (row swap)
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D8") "2.445"
Set-TextValue $ws.Range("E8") "-0.56%"

# Row 9 - LiechtensteinCryptoassetsExchange
Set-TextValue $ws.Range("D9") "0.1284"
Set-TextValue $ws.Range("E9") "26.73%"

# Row 10 - WazirX
Set-TextValue $ws.Range("D10") "0.1803"
Set-TextValue $ws.Range("E10") "1.72%"

# Row 11 - MandalaExchangeToken
Set-TextValue $ws.Range("D11") "0.09101"
Set-TextValue $ws.Range("E11") "0.18%"

# Row 12 - BitrueCoin
Set-TextValue $ws.Range("D12") "0.04354"
Set-TextValue $ws.Range("E12") "0.94%"

# Row 13 - BitMartToken
Set-TextValue $ws.Range("D13") "0.1047"
Set-TextValue $ws.Range("E13") "-0.60%"

# Row 14 - BitForexToken
Set-TextValue $ws.Range("D14") "0.001250"
Set-TextValue $ws.Range("E14") "-0.53%"

# Row 15 - TigerCash
Set-TextValue $ws.Range("D15") "0.005816"
Set-TextValue $ws.Range("E15") "0.17%"

# Row 16 - LEO
Set-TextValue $ws.Range("E16") "-0.11%"

# Row 17 - GateToken
Set-TextValue $ws.Range("D17") "4.289"
Set-TextValue $ws.Range("E17") "0.60%"

# Row 18 - BitpandaEcosystemToken
Set-TextValue $ws.Range("D18") "0.3317"
Set-TextValue $ws.Range("E18") "1.46%"

# Row 19 - MCDex
Set-TextValue $ws.Range("D19") "6.903"
Set-TextValue $ws.Range("E19") "1.54%"

# Row 20 - ProBitToken
Set-TextValue $ws.Range("E20") "2.72%"

# Row 21 - ZBToken
Set-TextValue $ws.Range("E21") "0.52%"

# Row 22 - CoinExToken
Set-TextValue $ws.Range("D22") "0.04054"
Set-TextValue $ws.Range("E22") "-2.47%"

# Row 23 - BitKan
Set-TextValue $ws.Range("D23") "0.001270"
Set-TextValue $ws.Range("E23") "3.69%"

# Row 24 - HotbitToken
Set-TextValue $ws.Range("D24") "0.004089"
Set-TextValue $ws.Range("E24") "-0.13%"

# Row 25 - NitroEx
Set-TextValue $ws.Range("D25") "0.0001272"
Set-TextValue $ws.Range("E25") "-2.21%"

# Row 26 - UpBots
Set-TextValue $ws.Range("E26") "24.69%"

# Row 38 - One
Set-TextValue $ws.Range("D38") "0.02418"
Set-TextValue $ws.Range("E38") "0.62%"

# Row 39 - IDEX
Set-TextValue $ws.Range("D39") "0.05230"
Set-TextValue $ws.Range("E39") "0.99%"

# Row 40 - KickToken
Set-TextValue $ws.Range("D40") "0.007850"
Set-TextValue $ws.Range("E40") "1.00%"

# Row 41 - BKEXToken
Set-TextValue $ws.Range("D41") "0.1302"
Set-TextValue $ws.Range("E41") "-0.54%"

# Row 42 - Dexo
Set-TextValue $ws.Range("D42") "0.006810"
Set-TextValue $ws.Range("E42") "-3.81%"

# Row 43 - CEJI
Set-TextValue $ws.Range("D43") "0.001863"
Set-TextValue $ws.Range("E43") "-2.86%"

# Row 44 - LocalTraders
Set-TextValue $ws.Range("D44") "0.007424"
Set-TextValue $ws.Range("E44") "-0.65%"

# Row 45 - PooCoin
Set-TextValue $ws.Range("D45") "0.3342"
Set-TextValue $ws.Range("E45") "-0.16%"

# Row 46 - CoinLion
Set-TextValue $ws.Range("D46") "0.00006874"
Set-TextValue $ws.Range("E46") "8.02%"

# Row 47 - Kangarootoken
Set-TextValue $ws.Range("E47") "0.09%"

# Row 48 - BOLO
Set-TextValue $ws.Range("D48") "0.1525"
Set-TextValue $ws.Range("E48") "2,108.87%"

# Row 49 - CoinbaseStockToken
Set-TextValue $ws.Range("E49") "-31.79%"

# Row 50 - CryptobidCoin
Set-TextValue $ws.Range("D50") "0.00002103"
Set-TextValue $ws.Range("E50") "0.09%"

# Row 51 - SpecialPowerGold
Set-TextValue $ws.Range("E51") "0.09%"
